$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Touch row 2 so it materializes as an (empty) row in the sheet,
#     matching the author re-saving with row 2 present but blank.
#     (No value is written - this keeps the row itself content-free.)
$ws.Cells.Item(2, 1).Font.Bold = $false

# --- Row 237 previously had placeholder (empty) entries in H/I;
#     the ticket was updated and those placeholders were cleared out.
$ws.Cells.Item(237, 8).ClearContents()
$ws.Cells.Item(237, 9).ClearContents()

# --- Append the four new ticket rows (238-241).
# Row 238
$ws.Cells.Item(238, 1).Value = "'2024-05-23"
$ws.Cells.Item(238, 2).Value = "15:21:03"
$ws.Cells.Item(238, 3).Value = "Palet atascado en la curva"
$ws.Cells.Item(238, 4).Value = "-"
$ws.Cells.Item(238, 5).Value = "-"
$ws.Cells.Item(238, 6).Value = "-"
$ws.Cells.Item(238, 7).Value = "-"

# Row 239
$ws.Cells.Item(239, 1).Value = "'2024-05-23"
$ws.Cells.Item(239, 2).Value = "15:33:56"
$ws.Cells.Item(239, 3).Value = "Palet atascado en la curva"
$ws.Cells.Item(239, 4).Value = "-"
$ws.Cells.Item(239, 5).Value = "-"
$ws.Cells.Item(239, 6).Value = "-"
$ws.Cells.Item(239, 7).Value = "-"
$ws.Cells.Item(239, 8).Value = "15:33:59"
$ws.Cells.Item(239, 9).Value = "0:00:03"

# Row 240
$ws.Cells.Item(240, 1).Value = "'2024-05-23"
$ws.Cells.Item(240, 2).Value = "15:34:01"
$ws.Cells.Item(240, 3).Value = "Fallo atornillador"
$ws.Cells.Item(240, 4).Value = "-"
$ws.Cells.Item(240, 5).Value = "-"
$ws.Cells.Item(240, 6).Value = "-"
$ws.Cells.Item(240, 7).Value = "-"

# Row 241
$ws.Cells.Item(241, 1).Value = "'2024-05-23"
$ws.Cells.Item(241, 2).Value = "15:36:48"
$ws.Cells.Item(241, 3).Value = "-"
$ws.Cells.Item(241, 4).Value = "-"
$ws.Cells.Item(241, 5).Value = "-"
$ws.Cells.Item(241, 6).Value = "Traza"
$ws.Cells.Item(241, 7).Value = "-"
$ws.Cells.Item(241, 8).Value = "15:36:49"
$ws.Cells.Item(241, 9).Value = "0:00:01"
